# Update Tff2-Cxcr4 NATMI LR-pair sheet with recomputed TPM values.
# Adds "MuSCs" as a new Sending cluster (with its own 4 rows, rows 10-13)
# and shifts the former "Resolving-Mac" sending-cluster rows down to rows
# 14-17 with freshly recomputed values, matching the updated TPM pipeline
# output. Existing rows 2-9 (ECs / FAPs senders) also get updated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Cells.Item(2,1).Value2 = "ECs"
$ws.Cells.Item(2,2).Value2 = "Tff2"
$ws.Cells.Item(2,3).Value2 = "Cxcr4"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 1
$ws.Cells.Item(2,6).Value2 = 0.3333333333333333
$ws.Cells.Item(2,7).Value2 = 0.151577
$ws.Cells.Item(2,8).Value2 = 0.454731
$ws.Cells.Item(2,9).Value2 = 0.1371818665598328
$ws.Cells.Item(2,10).Value2 = 0.1371818665598328
$ws.Cells.Item(2,11).Value2 = 3
$ws.Cells.Item(2,12).Value2 = 1
$ws.Cells.Item(2,13).Value2 = 25.69910333333333
$ws.Cells.Item(2,14).Value2 = 77.09731
$ws.Cells.Item(2,15).Value2 = 0.08761243344445813
$ws.Cells.Item(2,16).Value2 = 0.08761243344445814
$ws.Cells.Item(2,17).Value2 = 3.895392985956666
$ws.Cells.Item(2,18).Value2 = 35.05853687361
$ws.Cells.Item(2,19).Value2 = 0.01201883715375989
$ws.Cells.Item(2,20).Value2 = 0.01201883715375989

# Row 3
$ws.Cells.Item(3,1).Value2 = "ECs"
$ws.Cells.Item(3,2).Value2 = "Tff2"
$ws.Cells.Item(3,3).Value2 = "Cxcr4"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(3,5).Value2 = 1
$ws.Cells.Item(3,6).Value2 = 0.3333333333333333
$ws.Cells.Item(3,7).Value2 = 0.151577
$ws.Cells.Item(3,8).Value2 = 0.454731
$ws.Cells.Item(3,9).Value2 = 0.1371818665598328
$ws.Cells.Item(3,10).Value2 = 0.1371818665598328
$ws.Cells.Item(3,11).Value2 = 1
$ws.Cells.Item(3,12).Value2 = 0.3333333333333333
$ws.Cells.Item(3,13).Value2 = 0.1622346666666667
$ws.Cells.Item(3,14).Value2 = 0.486704
$ws.Cells.Item(3,15).Value2 = 0.0005530844306649811
$ws.Cells.Item(3,16).Value2 = 0.0005530844306649812
$ws.Cells.Item(3,17).Value2 = 0.02459104406933333
$ws.Cells.Item(3,18).Value2 = 0.221319396624
$ws.Cells.Item(3,19).Value2 = 0.00007587315456380454
$ws.Cells.Item(3,20).Value2 = 0.00007587315456380456

# Row 4
$ws.Cells.Item(4,1).Value2 = "ECs"
$ws.Cells.Item(4,2).Value2 = "Tff2"
$ws.Cells.Item(4,3).Value2 = "Cxcr4"
$ws.Cells.Item(4,4).Value2 = "MuSCs"
$ws.Cells.Item(4,5).Value2 = 1
$ws.Cells.Item(4,6).Value2 = 0.3333333333333333
$ws.Cells.Item(4,7).Value2 = 0.151577
$ws.Cells.Item(4,8).Value2 = 0.454731
$ws.Cells.Item(4,9).Value2 = 0.1371818665598328
$ws.Cells.Item(4,10).Value2 = 0.1371818665598328
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 46.92720933333334
$ws.Cells.Item(4,14).Value2 = 140.781628
$ws.Cells.Item(4,15).Value2 = 0.1599825079935015
$ws.Cells.Item(4,16).Value2 = 0.1599825079935016
$ws.Cells.Item(4,17).Value2 = 7.113085609118666
$ws.Cells.Item(4,18).Value2 = 64.017770482068
$ws.Cells.Item(4,19).Value2 = 0.02194669906347191
$ws.Cells.Item(4,20).Value2 = 0.02194669906347192

# Row 5
$ws.Cells.Item(5,1).Value2 = "ECs"
$ws.Cells.Item(5,2).Value2 = "Tff2"
$ws.Cells.Item(5,3).Value2 = "Cxcr4"
$ws.Cells.Item(5,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(5,5).Value2 = 1
$ws.Cells.Item(5,6).Value2 = 0.3333333333333333
$ws.Cells.Item(5,7).Value2 = 0.151577
$ws.Cells.Item(5,8).Value2 = 0.454731
$ws.Cells.Item(5,9).Value2 = 0.1371818665598328
$ws.Cells.Item(5,10).Value2 = 0.1371818665598328
$ws.Cells.Item(5,11).Value2 = 3
$ws.Cells.Item(5,12).Value2 = 1
$ws.Cells.Item(5,13).Value2 = 220.538579
$ws.Cells.Item(5,14).Value2 = 661.615737
$ws.Cells.Item(5,15).Value2 = 0.7518519741313753
$ws.Cells.Item(5,16).Value2 = 0.7518519741313754
$ws.Cells.Item(5,17).Value2 = 33.428576189083
$ws.Cells.Item(5,18).Value2 = 300.857185701747
$ws.Cells.Item(5,19).Value2 = 0.1031404571880372
$ws.Cells.Item(5,20).Value2 = 0.1031404571880372

# Row 6
$ws.Cells.Item(6,1).Value2 = "FAPs"
$ws.Cells.Item(6,2).Value2 = "Tff2"
$ws.Cells.Item(6,3).Value2 = "Cxcr4"
$ws.Cells.Item(6,4).Value2 = "ECs"
$ws.Cells.Item(6,5).Value2 = 1
$ws.Cells.Item(6,6).Value2 = 0.3333333333333333
$ws.Cells.Item(6,7).Value2 = 0.045591
$ws.Cells.Item(6,8).Value2 = 0.136773
$ws.Cells.Item(6,9).Value2 = 0.04126126310937238
$ws.Cells.Item(6,10).Value2 = 0.04126126310937238
$ws.Cells.Item(6,11).Value2 = 3
$ws.Cells.Item(6,12).Value2 = 1
$ws.Cells.Item(6,13).Value2 = 25.69910333333333
$ws.Cells.Item(6,14).Value2 = 77.09731
$ws.Cells.Item(6,15).Value2 = 0.08761243344445813
$ws.Cells.Item(6,16).Value2 = 0.08761243344445814
$ws.Cells.Item(6,17).Value2 = 1.17164782007
$ws.Cells.Item(6,18).Value2 = 10.54483038063
$ws.Cells.Item(6,19).Value2 = 0.003614999668004163
$ws.Cells.Item(6,20).Value2 = 0.003614999668004164

# Row 7
$ws.Cells.Item(7,1).Value2 = "FAPs"
$ws.Cells.Item(7,2).Value2 = "Tff2"
$ws.Cells.Item(7,3).Value2 = "Cxcr4"
$ws.Cells.Item(7,4).Value2 = "FAPs"
$ws.Cells.Item(7,5).Value2 = 1
$ws.Cells.Item(7,6).Value2 = 0.3333333333333333
$ws.Cells.Item(7,7).Value2 = 0.045591
$ws.Cells.Item(7,8).Value2 = 0.136773
$ws.Cells.Item(7,9).Value2 = 0.04126126310937238
$ws.Cells.Item(7,10).Value2 = 0.04126126310937238
$ws.Cells.Item(7,11).Value2 = 1
$ws.Cells.Item(7,12).Value2 = 0.3333333333333333
$ws.Cells.Item(7,13).Value2 = 0.1622346666666667
$ws.Cells.Item(7,14).Value2 = 0.486704
$ws.Cells.Item(7,15).Value2 = 0.0005530844306649811
$ws.Cells.Item(7,16).Value2 = 0.0005530844306649812
$ws.Cells.Item(7,17).Value2 = 0.007396440688
$ws.Cells.Item(7,18).Value2 = 0.066567966192
$ws.Cells.Item(7,19).Value2 = 0.00002282096221536521
$ws.Cells.Item(7,20).Value2 = 0.00002282096221536521

# Row 8
$ws.Cells.Item(8,1).Value2 = "FAPs"
$ws.Cells.Item(8,2).Value2 = "Tff2"
$ws.Cells.Item(8,3).Value2 = "Cxcr4"
$ws.Cells.Item(8,4).Value2 = "MuSCs"
$ws.Cells.Item(8,5).Value2 = 1
$ws.Cells.Item(8,6).Value2 = 0.3333333333333333
$ws.Cells.Item(8,7).Value2 = 0.045591
$ws.Cells.Item(8,8).Value2 = 0.136773
$ws.Cells.Item(8,9).Value2 = 0.04126126310937238
$ws.Cells.Item(8,10).Value2 = 0.04126126310937238
$ws.Cells.Item(8,11).Value2 = 3
$ws.Cells.Item(8,12).Value2 = 1
$ws.Cells.Item(8,13).Value2 = 46.92720933333334
$ws.Cells.Item(8,14).Value2 = 140.781628
$ws.Cells.Item(8,15).Value2 = 0.1599825079935015
$ws.Cells.Item(8,16).Value2 = 0.1599825079935016
$ws.Cells.Item(8,17).Value2 = 2.139458400716
$ws.Cells.Item(8,18).Value2 = 19.255125606444
$ws.Cells.Item(8,19).Value2 = 0.006601080355217137
$ws.Cells.Item(8,20).Value2 = 0.006601080355217137

# Row 9
$ws.Cells.Item(9,1).Value2 = "FAPs"
$ws.Cells.Item(9,2).Value2 = "Tff2"
$ws.Cells.Item(9,3).Value2 = "Cxcr4"
$ws.Cells.Item(9,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(9,5).Value2 = 1
$ws.Cells.Item(9,6).Value2 = 0.3333333333333333
$ws.Cells.Item(9,7).Value2 = 0.045591
$ws.Cells.Item(9,8).Value2 = 0.136773
$ws.Cells.Item(9,9).Value2 = 0.04126126310937238
$ws.Cells.Item(9,10).Value2 = 0.04126126310937238
$ws.Cells.Item(9,11).Value2 = 3
$ws.Cells.Item(9,12).Value2 = 1
$ws.Cells.Item(9,13).Value2 = 220.538579
$ws.Cells.Item(9,14).Value2 = 661.615737
$ws.Cells.Item(9,15).Value2 = 0.7518519741313753
$ws.Cells.Item(9,16).Value2 = 0.7518519741313754
$ws.Cells.Item(9,17).Value2 = 10.054574355189
$ws.Cells.Item(9,18).Value2 = 90.491169196701
$ws.Cells.Item(9,19).Value2 = 0.03102236212393571
$ws.Cells.Item(9,20).Value2 = 0.03102236212393571

# Row 10
$ws.Cells.Item(10,1).Value2 = "MuSCs"
$ws.Cells.Item(10,2).Value2 = "Tff2"
$ws.Cells.Item(10,3).Value2 = "Cxcr4"
$ws.Cells.Item(10,4).Value2 = "ECs"
$ws.Cells.Item(10,5).Value2 = 2
$ws.Cells.Item(10,6).Value2 = 0.6666666666666666
$ws.Cells.Item(10,7).Value2 = 0.1599346666666667
$ws.Cells.Item(10,8).Value2 = 0.479804
$ws.Cells.Item(10,9).Value2 = 0.1447458130254459
$ws.Cells.Item(10,10).Value2 = 0.1447458130254458
$ws.Cells.Item(10,11).Value2 = 3
$ws.Cells.Item(10,12).Value2 = 1
$ws.Cells.Item(10,13).Value2 = 25.69910333333333
$ws.Cells.Item(10,14).Value2 = 77.09731
$ws.Cells.Item(10,15).Value2 = 0.08761243344445813
$ws.Cells.Item(10,16).Value2 = 0.08761243344445814
$ws.Cells.Item(10,17).Value2 = 4.110177525248888
$ws.Cells.Item(10,18).Value2 = 36.99159772724
$ws.Cells.Item(10,19).Value2 = 0.01268153291005586
$ws.Cells.Item(10,20).Value2 = 0.01268153291005586

# Row 11
$ws.Cells.Item(11,1).Value2 = "MuSCs"
$ws.Cells.Item(11,2).Value2 = "Tff2"
$ws.Cells.Item(11,3).Value2 = "Cxcr4"
$ws.Cells.Item(11,4).Value2 = "FAPs"
$ws.Cells.Item(11,5).Value2 = 2
$ws.Cells.Item(11,6).Value2 = 0.6666666666666666
$ws.Cells.Item(11,7).Value2 = 0.1599346666666667
$ws.Cells.Item(11,8).Value2 = 0.479804
$ws.Cells.Item(11,9).Value2 = 0.1447458130254459
$ws.Cells.Item(11,10).Value2 = 0.1447458130254458
$ws.Cells.Item(11,11).Value2 = 1
$ws.Cells.Item(11,12).Value2 = 0.3333333333333333
$ws.Cells.Item(11,13).Value2 = 0.1622346666666667
$ws.Cells.Item(11,14).Value2 = 0.486704
$ws.Cells.Item(11,15).Value2 = 0.0005530844306649811
$ws.Cells.Item(11,16).Value2 = 0.0005530844306649812
$ws.Cells.Item(11,17).Value2 = 0.02594694733511111
$ws.Cells.Item(11,18).Value2 = 0.233522526016
$ws.Cells.Item(11,19).Value2 = 0.00008005665558831854
$ws.Cells.Item(11,20).Value2 = 0.00008005665558831854

# Row 12
$ws.Cells.Item(12,1).Value2 = "MuSCs"
$ws.Cells.Item(12,2).Value2 = "Tff2"
$ws.Cells.Item(12,3).Value2 = "Cxcr4"
$ws.Cells.Item(12,4).Value2 = "MuSCs"
$ws.Cells.Item(12,5).Value2 = 2
$ws.Cells.Item(12,6).Value2 = 0.6666666666666666
$ws.Cells.Item(12,7).Value2 = 0.1599346666666667
$ws.Cells.Item(12,8).Value2 = 0.479804
$ws.Cells.Item(12,9).Value2 = 0.1447458130254459
$ws.Cells.Item(12,10).Value2 = 0.1447458130254458
$ws.Cells.Item(12,11).Value2 = 3
$ws.Cells.Item(12,12).Value2 = 1
$ws.Cells.Item(12,13).Value2 = 46.92720933333334
$ws.Cells.Item(12,14).Value2 = 140.781628
$ws.Cells.Item(12,15).Value2 = 0.1599825079935015
$ws.Cells.Item(12,16).Value2 = 0.1599825079935016
$ws.Cells.Item(12,17).Value2 = 7.505287582323557
$ws.Cells.Item(12,18).Value2 = 67.547588240912
$ws.Cells.Item(12,19).Value2 = 0.02315679818936927
$ws.Cells.Item(12,20).Value2 = 0.02315679818936927

# Row 13
$ws.Cells.Item(13,1).Value2 = "MuSCs"
$ws.Cells.Item(13,2).Value2 = "Tff2"
$ws.Cells.Item(13,3).Value2 = "Cxcr4"
$ws.Cells.Item(13,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(13,5).Value2 = 2
$ws.Cells.Item(13,6).Value2 = 0.6666666666666666
$ws.Cells.Item(13,7).Value2 = 0.1599346666666667
$ws.Cells.Item(13,8).Value2 = 0.479804
$ws.Cells.Item(13,9).Value2 = 0.1447458130254459
$ws.Cells.Item(13,10).Value2 = 0.1447458130254458
$ws.Cells.Item(13,11).Value2 = 3
$ws.Cells.Item(13,12).Value2 = 1
$ws.Cells.Item(13,13).Value2 = 220.538579
$ws.Cells.Item(13,14).Value2 = 661.615737
$ws.Cells.Item(13,15).Value2 = 0.7518519741313753
$ws.Cells.Item(13,16).Value2 = 0.7518519741313754
$ws.Cells.Item(13,17).Value2 = 35.27176411950533
$ws.Cells.Item(13,18).Value2 = 317.445877075548
$ws.Cells.Item(13,19).Value2 = 0.1088274252704324
$ws.Cells.Item(13,20).Value2 = 0.1088274252704324

# Row 14
$ws.Cells.Item(14,1).Value2 = "Resolving-Mac"
$ws.Cells.Item(14,2).Value2 = "Tff2"
$ws.Cells.Item(14,3).Value2 = "Cxcr4"
$ws.Cells.Item(14,4).Value2 = "ECs"
$ws.Cells.Item(14,5).Value2 = 2
$ws.Cells.Item(14,6).Value2 = 0.6666666666666666
$ws.Cells.Item(14,7).Value2 = 0.7478319999999999
$ws.Cells.Item(14,8).Value2 = 2.243496
$ws.Cells.Item(14,9).Value2 = 0.676811057305349
$ws.Cells.Item(14,10).Value2 = 0.676811057305349
$ws.Cells.Item(14,11).Value2 = 3
$ws.Cells.Item(14,12).Value2 = 1
$ws.Cells.Item(14,13).Value2 = 25.69910333333333
$ws.Cells.Item(14,14).Value2 = 77.09731
$ws.Cells.Item(14,15).Value2 = 0.08761243344445813
$ws.Cells.Item(14,16).Value2 = 0.08761243344445814
$ws.Cells.Item(14,17).Value2 = 19.21861184397333
$ws.Cells.Item(14,18).Value2 = 172.96750659576
$ws.Cells.Item(14,19).Value2 = 0.05929706371263822
$ws.Cells.Item(14,20).Value2 = 0.05929706371263823

# Row 15
$ws.Cells.Item(15,1).Value2 = "Resolving-Mac"
$ws.Cells.Item(15,2).Value2 = "Tff2"
$ws.Cells.Item(15,3).Value2 = "Cxcr4"
$ws.Cells.Item(15,4).Value2 = "FAPs"
$ws.Cells.Item(15,5).Value2 = 2
$ws.Cells.Item(15,6).Value2 = 0.6666666666666666
$ws.Cells.Item(15,7).Value2 = 0.7478319999999999
$ws.Cells.Item(15,8).Value2 = 2.243496
$ws.Cells.Item(15,9).Value2 = 0.676811057305349
$ws.Cells.Item(15,10).Value2 = 0.676811057305349
$ws.Cells.Item(15,11).Value2 = 1
$ws.Cells.Item(15,12).Value2 = 0.3333333333333333
$ws.Cells.Item(15,13).Value2 = 0.1622346666666667
$ws.Cells.Item(15,14).Value2 = 0.486704
$ws.Cells.Item(15,15).Value2 = 0.0005530844306649811
$ws.Cells.Item(15,16).Value2 = 0.0005530844306649812
$ws.Cells.Item(15,17).Value2 = 0.1213242752426667
$ws.Cells.Item(15,18).Value2 = 1.091918477184
$ws.Cells.Item(15,19).Value2 = 0.0003743336582974928
$ws.Cells.Item(15,20).Value2 = 0.0003743336582974929

# Row 16
$ws.Cells.Item(16,1).Value2 = "Resolving-Mac"
$ws.Cells.Item(16,2).Value2 = "Tff2"
$ws.Cells.Item(16,3).Value2 = "Cxcr4"
$ws.Cells.Item(16,4).Value2 = "MuSCs"
$ws.Cells.Item(16,5).Value2 = 2
$ws.Cells.Item(16,6).Value2 = 0.6666666666666666
$ws.Cells.Item(16,7).Value2 = 0.7478319999999999
$ws.Cells.Item(16,8).Value2 = 2.243496
$ws.Cells.Item(16,9).Value2 = 0.676811057305349
$ws.Cells.Item(16,10).Value2 = 0.676811057305349
$ws.Cells.Item(16,11).Value2 = 3
$ws.Cells.Item(16,12).Value2 = 1
$ws.Cells.Item(16,13).Value2 = 46.92720933333334
$ws.Cells.Item(16,14).Value2 = 140.781628
$ws.Cells.Item(16,15).Value2 = 0.1599825079935015
$ws.Cells.Item(16,16).Value2 = 0.1599825079935016
$ws.Cells.Item(16,17).Value2 = 35.09366881016533
$ws.Cells.Item(16,18).Value2 = 315.843019291488
$ws.Cells.Item(16,19).Value2 = 0.1082779303854432
$ws.Cells.Item(16,20).Value2 = 0.1082779303854432

# Row 17
$ws.Cells.Item(17,1).Value2 = "Resolving-Mac"
$ws.Cells.Item(17,2).Value2 = "Tff2"
$ws.Cells.Item(17,3).Value2 = "Cxcr4"
$ws.Cells.Item(17,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(17,5).Value2 = 2
$ws.Cells.Item(17,6).Value2 = 0.6666666666666666
$ws.Cells.Item(17,7).Value2 = 0.7478319999999999
$ws.Cells.Item(17,8).Value2 = 2.243496
$ws.Cells.Item(17,9).Value2 = 0.676811057305349
$ws.Cells.Item(17,10).Value2 = 0.676811057305349
$ws.Cells.Item(17,11).Value2 = 3
$ws.Cells.Item(17,12).Value2 = 1
$ws.Cells.Item(17,13).Value2 = 220.538579
$ws.Cells.Item(17,14).Value2 = 661.615737
$ws.Cells.Item(17,15).Value2 = 0.7518519741313753
$ws.Cells.Item(17,16).Value2 = 0.7518519741313754
$ws.Cells.Item(17,17).Value2 = 164.925806610728
$ws.Cells.Item(17,18).Value2 = 1484.332259496552
$ws.Cells.Item(17,19).Value2 = 0.50886172954897
$ws.Cells.Item(17,20).Value2 = 0.50886172954897
